$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p069r_2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p069r_2</id>", 2)
